$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern used throughout: the Price (D) column stores plain-text
# numeric-looking strings (e.g. "1.001", "22.24"). A bare `.Value = "..."`
# assignment lets Excel reinterpret anything that parses cleanly as a number
# (General number format), which would silently convert the cell from Text
# to Number and mangle values like "1.000" -> 1 or "2.180" -> 2.18. To keep
# every Price cell text-typed (matching the source workbook), each write
# temporarily forces NumberFormat to Text ("@"), assigns the literal string,
# then resets Style back to "Normal" so no stray number-format style lingers.

# --- Row 12 / Row 13: BinanceUSD and Solana swap places ---
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.31%  "

$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.000"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.14%  "

# --- Remaining per-row Price (D) / Volume(1h) (E) refreshes ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.363.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4451"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3723"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").Value = "  +4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.138"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.334"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.636"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.853.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.96%  "
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06532"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  +4.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.257"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.408.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.65%  "
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.180"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.049.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.325"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.213"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.945"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.63%  "
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.634"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("E36").Value = "  +4.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2201"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.215"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6626"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("E40").Value = "  +3.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.203"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.191"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.439"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6162"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.70%  "
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.049"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.163"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07013"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.27%  "
